$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws.Range("H17").Value = 844917.2
$ws.Range("I17").Value = 410.77777
$ws.Range("J17").Value = 1231386.2
$ws.Range("K17").Value = 1232.33331
$ws.Range("L17").Value = 3694158.6
$ws.Range("M17").Value = -1064.33331
$ws.Range("N17").Value = -3694494.6

$ws.Range("H42").Value = 409.4
$ws.Range("I42").Value = 74
$ws.Range("J42").Value = 633
$ws.Range("K42").Value = 222
$ws.Range("L42").Value = 1899
$ws.Range("M42").Value = 8
$ws.Range("N42").Value = -2359

$ws.Range("H135").Value = 852.31915
$ws.Range("I135").Value = 475.07895
$ws.Range("J135").Value = 2445.111
$ws.Range("K135").Value = 4275.71055
$ws.Range("L135").Value = 22005.999
$ws.Range("M135").Value = -1740.71055
$ws.Range("N135").Value = -27075.999

$ws.Range("H137").Value = 1262.1428
$ws.Range("I137").Value = 1225.75
$ws.Range("J137").Value = 1300.6765
$ws.Range("K137").Value = 3677.25
$ws.Range("L137").Value = 3902.0295
$ws.Range("M137").Value = -1127.25
$ws.Range("N137").Value = -9002.029500000001

$ws.Range("H138").Value = 1070.35
$ws.Range("I138").Value = 562.60785
$ws.Range("J138").Value = 1598.8163
$ws.Range("K138").Value = 1687.82355
$ws.Range("L138").Value = 4796.448899999999
$ws.Range("M138").Value = 3452.17645
$ws.Range("N138").Value = -15076.4489

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1327.409
$ws.Range("I2").Value = 1390.5555
$ws.Range("K2").Value = 1390.5555
$ws.Range("M2").Value = -1277.5555

$ws.Range("H10").Value = 60002.5
$ws.Range("I10").Value = 50000
$ws.Range("J10").Value = 70005
$ws.Range("K10").Value = 50000
$ws.Range("L10").Value = 70005
$ws.Range("M10").Value = -49830
$ws.Range("N10").Value = -70345

$ws.Range("H32").Value = 4580.83
$ws.Range("I32").Value = 4071.7866
$ws.Range("J32").Value = 8699.454
$ws.Range("K32").Value = 4071.7866
$ws.Range("L32").Value = 8699.454
$ws.Range("M32").Value = -3784.7866
$ws.Range("N32").Value = -9273.454

$ws.Range("H53").Value = 3000
$ws.Range("I53").Value = 3000
$ws.Range("K53").Value = 3000
$ws.Range("M53").Value = -2318

$ws.Range("H88").Value = 1462.375
$ws.Range("I88").Value = 1363.4546
$ws.Range("J88").Value = 1680
$ws.Range("K88").Value = 1363.4546
$ws.Range("L88").Value = 1680
$ws.Range("M88").Value = -957.4546
$ws.Range("N88").Value = -2492

$ws.Range("H91").Value = 1462.375
$ws.Range("I91").Value = 1363.4546
$ws.Range("J91").Value = 1680
$ws.Range("K91").Value = 1363.4546
$ws.Range("L91").Value = 1680
$ws.Range("M91").Value = 40.54539999999997
$ws.Range("N91").Value = -4488

$ws.Range("H97").Value = 1102.6552
$ws.Range("I97").Value = 1023.3684
$ws.Range("J97").Value = 1253.3
$ws.Range("K97").Value = 1023.3684
$ws.Range("L97").Value = 1253.3
$ws.Range("M97").Value = -527.3684
$ws.Range("N97").Value = -2245.3

$ws.Range("H115").Value = 65000
$ws.Range("J115").Value = 65000
$ws.Range("L115").Value = 65000
$ws.Range("N115").Value = -68134

$ws.Range("H116").Value = 1327.409
$ws.Range("I116").Value = 1390.5555
$ws.Range("K116").Value = 1390.5555
$ws.Range("M116").Value = 903.4445000000001

$ws.Range("H122").Value = 169437.67
$ws.Range("I122").Value = 501456
$ws.Range("J122").Value = 3428.5
$ws.Range("K122").Value = 1504368
$ws.Range("L122").Value = 10285.5
$ws.Range("M122").Value = -1501918
$ws.Range("N122").Value = -15185.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1327.409
$ws.Range("I3").Value = 1390.5555
$ws.Range("K3").Value = 1390.5555
$ws.Range("M3").Value = -1276.5555

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 26500.666
$ws.Range("J4").Value = 26500.666
$ws.Range("L4").Value = 26500.666
$ws.Range("N4").Value = -26724.666

$ws.Range("H31").Value = 4377.875
$ws.Range("I31").Value = 1331.0938
$ws.Range("J31").Value = 7424.6562
$ws.Range("K31").Value = 1331.0938
$ws.Range("L31").Value = 7424.6562
$ws.Range("M31").Value = -1036.0938
$ws.Range("N31").Value = -8014.6562

$ws.Range("H34").Value = 4377.875
$ws.Range("I34").Value = 1331.0938
$ws.Range("J34").Value = 7424.6562
$ws.Range("K34").Value = 1331.0938
$ws.Range("L34").Value = 7424.6562
$ws.Range("M34").Value = -1129.0938
$ws.Range("N34").Value = -7828.6562

$ws.Range("H58").Value = 1028.1346
$ws.Range("I58").Value = 768.6667
$ws.Range("K58").Value = 768.6667
$ws.Range("M58").Value = -565.6667

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H134").Value = 3665.5107
$ws.Range("I134").Value = 3669.8684
$ws.Range("J134").Value = 3647.111
$ws.Range("K134").Value = 11009.6052
$ws.Range("L134").Value = 10941.333
$ws.Range("M134").Value = -8474.6052
$ws.Range("N134").Value = -16011.333

$ws.Range("H136").Value = 1028.1346
$ws.Range("I136").Value = 768.6667
$ws.Range("K136").Value = 2306.0001
$ws.Range("M136").Value = 243.9998999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7779277.5
$ws.Range("I4").Value = 4001000
$ws.Range("J4").Value = 12502125
$ws.Range("K4").Value = 12003000
$ws.Range("L4").Value = 37506375
$ws.Range("M4").Value = -12002888
$ws.Range("N4").Value = -37506599

$ws.Range("H5").Value = 1296.8857
$ws.Range("I5").Value = 394.55
$ws.Range("J5").Value = 2500
$ws.Range("K5").Value = 1183.65
$ws.Range("L5").Value = 7500
$ws.Range("M5").Value = -1071.65
$ws.Range("N5").Value = -7724

$ws.Range("H63").Value = 945
$ws.Range("I63").Value = 945
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2835
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2086
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 945
$ws.Range("I66").Value = 945
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 8505
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -4761
$ws.Range("N66").ClearContents()

$ws.Range("H75").Value = 5171.6665
$ws.Range("J75").Value = 7657.5
$ws.Range("L75").Value = 22972.5
$ws.Range("N75").Value = -24968.5

$ws.Range("H78").Value = 5171.6665
$ws.Range("J78").Value = 7657.5
$ws.Range("L78").Value = 68917.5
$ws.Range("N78").Value = -78901.5

$ws.Range("H102").Value = 1509.6666
$ws.Range("J102").Value = 1014.5
$ws.Range("L102").Value = 3043.5
$ws.Range("N102").Value = -7911.5

$ws.Range("H103").Value = 1921.2
$ws.Range("J103").Value = 3018.6667
$ws.Range("L103").Value = 9056.000100000001
$ws.Range("N103").Value = -10814.0001

$ws.Range("H107").Value = 20000304
$ws.Range("I107").Value = 286.78125
$ws.Range("J107").Value = 55555892
$ws.Range("K107").Value = 860.34375
$ws.Range("L107").Value = 166667676
$ws.Range("M107").Value = 1059.65625
$ws.Range("N107").Value = -166671516

$ws.Range("H112").Value = 4286.3125
$ws.Range("I112").Value = 2642.3333
$ws.Range("J112").Value = 6400
$ws.Range("K112").Value = 7926.999899999999
$ws.Range("L112").Value = 19200
$ws.Range("M112").Value = -6818.999899999999
$ws.Range("N112").Value = -21416

$ws.Range("H122").Value = 2710.681
$ws.Range("I122").Value = 378.7
$ws.Range("J122").Value = 6825.9414
$ws.Range("K122").Value = 3408.3
$ws.Range("L122").Value = 61433.47259999999
$ws.Range("M122").Value = -958.2999999999997
$ws.Range("N122").Value = -66333.47259999999

$ws.Range("H131").Value = 2965.9312
$ws.Range("I131").Value = 413.42856
$ws.Range("J131").Value = 3778.0908
$ws.Range("K131").Value = 1240.28568
$ws.Range("L131").Value = 11334.2724
$ws.Range("M131").Value = 3799.71432
$ws.Range("N131").Value = -21414.2724

$ws.Range("H135").Value = 1296.8857
$ws.Range("I135").Value = 394.55
$ws.Range("J135").Value = 2500
$ws.Range("K135").Value = 3550.95
$ws.Range("L135").Value = 22500
$ws.Range("M135").Value = -1015.95
$ws.Range("N135").Value = -27570

$ws.Range("H137").Value = 6179198.5
$ws.Range("I137").Value = 15161795
$ws.Range("J137").Value = 3663.375
$ws.Range("K137").Value = 45485385
$ws.Range("L137").Value = 10990.125
$ws.Range("M137").Value = -45480285
$ws.Range("N137").Value = -21190.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2000166.6
$ws.Range("I3").Value = 2500250
$ws.Range("J3").Value = 1000000
$ws.Range("K3").Value = 2500250
$ws.Range("L3").Value = 1000000
$ws.Range("M3").Value = -2500134
$ws.Range("N3").Value = -1000232

$ws.Range("H102").Value = 1794.027
$ws.Range("I102").Value = 1685.4828
$ws.Range("J102").Value = 2187.5
$ws.Range("K102").Value = 1685.4828
$ws.Range("L102").Value = 2187.5
$ws.Range("M102").Value = -63.4828
$ws.Range("N102").Value = -5431.5

$ws.Range("H132").Value = 2882.1333
$ws.Range("I132").Value = 2479.8215
$ws.Range("J132").Value = 3544.7646
$ws.Range("K132").Value = 7439.4645
$ws.Range("L132").Value = 10634.2938
$ws.Range("M132").Value = -4909.4645
$ws.Range("N132").Value = -15694.2938

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2408001.5
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 2408001.5
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 2408001.5
$ws.Range("N2").Value = -2408225.5
$ws.Range("M2").ClearContents()

$ws.Range("H13").Value = 70007
$ws.Range("J13").Value = 70007
$ws.Range("L13").Value = 70007
$ws.Range("N13").Value = -70287

$ws.Range("H100").Value = 3613.3333
$ws.Range("I100").Value = 3522.8572
$ws.Range("J100").Value = 3740
$ws.Range("K100").Value = 3522.8572
$ws.Range("L100").Value = 3740
$ws.Range("M100").Value = -2981.8572
$ws.Range("N100").Value = -4822

$ws.Range("H132").Value = 2363.3235
$ws.Range("I132").Value = 2190.3726
$ws.Range("K132").Value = 6571.1178
$ws.Range("M132").Value = -4041.1178

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 251250990
$ws.Range("J2").Value = 333334660
$ws.Range("L2").Value = 333334660
$ws.Range("N2").Value = -333334884

$ws.Range("H20").Value = 38915.07
$ws.Range("I20").Value = 400
$ws.Range("J20").Value = 45334.25
$ws.Range("K20").Value = 400
$ws.Range("L20").Value = 45334.25
$ws.Range("M20").Value = -160
$ws.Range("N20").Value = -45814.25
